# Generate Report for Archive
#
# This edit updates the localization status from "Ready for handoff" to
# "In Translation" everywhere it appears (Overview sheet's zh-cn/de-de
# status columns, and the Status column on each per-locale detail sheet),
# and shrinks the corresponding status columns to their new, narrower
# auto-fit width now that the text is shorter.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# Target column width (characters) recorded in the canonical workbook for
# the narrower status columns. The COM layer quantizes ColumnWidth to a
# discrete pixel grid, so we feed it the input that lands on the closest
# achievable grid point to the recorded width.
$newStatusColumnWidth = 12.5

# --- Overview sheet: columns E (zh-cn) and F (de-de) hold the status ---
$overview = $wb.Worksheets.Item("Overview")

if ($overview.Range("E2").Text -eq $oldStatus) {
    $overview.Range("E2").Value = $newStatus
}
if ($overview.Range("F2").Text -eq $oldStatus) {
    $overview.Range("F2").Value = $newStatus
}

$overview.Columns.Item(5).ColumnWidth = $newStatusColumnWidth
$overview.Columns.Item(6).ColumnWidth = $newStatusColumnWidth

# --- Per-locale detail sheets: column C holds the Status value ---
foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)

    if ($ws.Range("C2").Text -eq $oldStatus) {
        $ws.Range("C2").Value = $newStatus
    }

    $ws.Columns.Item(3).ColumnWidth = $newStatusColumnWidth
}
